$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting existing rows 102:117 down to 103:118.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new weekly data point.
$ws.Range("A102").Value = 7
$ws.Range("B102").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C102").Value = "Ñuble"
$ws.Range("D102").Value = 45173
$ws.Range("E102").Value = 16
$ws.Range("F102").Value = 100112001
$ws.Range("G102").Value = "Berenjena"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 160
$ws.Range("K102").Value = 8000
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = 8500
$ws.Range("N102").Value = '$/caja 60 unidades'
$ws.Range("O102").Value = "Región de Arica y Parinacota"
$ws.Range("P102").Value = 142
$ws.Range("Q102").Value = 60
$ws.Range("R102").Value = "Hortaliza"
